$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains values that look numeric (e.g. "216.79")
# as well as locale-formatted numbers with multiple separators
# (e.g. "26.106.99"). Force the column to Text format first so that
# Excel does not auto-convert the numeric-looking values to actual
# numbers when we assign the new string values below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.106.99'
$ws.Range("E2").Value = '  -1.26%  '
$ws.Range("D3").Value = '1.668.30'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '216.79'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = '0.5117'
$ws.Range("E6").Value = '  +3.49%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '0.2640'
$ws.Range("E8").Value = '  +1.68%  '
$ws.Range("D9").Value = '0.06429'
$ws.Range("E9").Value = '  +5.05%  '
$ws.Range("D10").Value = '21.64'
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("E11").Value = '  +1.94%  '
$ws.Range("D12").Value = '1.664.08'
$ws.Range("E12").Value = '  -0.45%  '
$ws.Range("D13").Value = '4.511'
$ws.Range("E13").Value = '  +2.25%  '
$ws.Range("D14").Value = '0.5817'
$ws.Range("E14").Value = '  +1.75%  '
$ws.Range("D15").Value = '0.000008583'
$ws.Range("E15").Value = '  +5.14%  '
$ws.Range("D16").Value = '64.44'
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("D17").Value = '26.165.14'
$ws.Range("E17").Value = '  -1.08%  '
$ws.Range("D18").Value = '4.932'
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("E20").Value = '  +1.40%  '
$ws.Range("D21").Value = '189.59'
$ws.Range("E21").Value = '  +3.97%  '
$ws.Range("D22").Value = '6.217'
$ws.Range("E22").Value = '  +1.11%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = '145.72'
$ws.Range("E24").Value = '  +1.07%  '
$ws.Range("D25").Value = '7.640'
$ws.Range("E25").Value = '  +1.91%  '
$ws.Range("D26").Value = '0.1205'
$ws.Range("E26").Value = '  +7.16%  '
$ws.Range("D27").Value = '15.61'
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("D28").Value = '0.06407'
$ws.Range("E28").Value = '  +13.53%  '
$ws.Range("D29").Value = '1.299'
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("D30").Value = '1.318'
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").Value = '3.525'
$ws.Range("E31").Value = '  +1.65%  '
$ws.Range("D32").Value = '3.514'
$ws.Range("E32").Value = '  +1.90%  '
$ws.Range("D33").Value = '1.634'
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("D34").Value = '1.019'
$ws.Range("E34").Value = '  +1.72%  '
$ws.Range("D35").Value = '0.6092'
$ws.Range("E35").Value = '  +3.99%  '
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("D37").Value = '2.649'
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("D38").Value = '6.187'
$ws.Range("E38").Value = '  +5.33%  '
$ws.Range("E39").Value = '  +2.03%  '
$ws.Range("D40").Value = '1.080.56'
$ws.Range("E40").Value = '  +1.01%  '
$ws.Range("D41").Value = '0.8606'
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("D42").Value = '1.009'
$ws.Range("E42").Value = '  +0.75%  '
$ws.Range("D43").Value = '101.04'
$ws.Range("E43").Value = '  +3.15%  '
$ws.Range("D44").Value = '1.815.27'
$ws.Range("E44").Value = '  -1.11%  '
$ws.Range("D45").Value = '0.00000000111'
$ws.Range("E45").Value = '  +6.62%  '
$ws.Range("D46").Value = '56.29'
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").Value = '1.011'
$ws.Range("E47").Value = '  +0.49%  '
$ws.Range("D48").Value = '8.072'
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").Value = '0.05204'
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("D50").Value = '0.4285'
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("D51").Value = '5.953'
$ws.Range("E51").Value = '  +6.79%  '
